# Updates the cryptocurrency price (column D) and 1h volume change (column E)
# figures on the active worksheet to reflect the latest scrape, matching the
# commit "Updated cryptos list on Thu Jun 29 23:51:32 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.465.92"
$ws.Range("E2").Value = "  +1.15%  "
$ws.Range("D3").Value = "'1.853.28"
$ws.Range("E3").Value = "  +1.20%  "
$ws.Range("D4").Value = "'0.9999"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'233.45"
$ws.Range("E5").Value = "  +1.15%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").Value = "'0.4746"
$ws.Range("E7").Value = "  +1.91%  "
$ws.Range("D8").Value = "'0.2757"
$ws.Range("E8").Value = "  +3.13%  "
$ws.Range("D9").Value = "'0.06346"
$ws.Range("E9").Value = "  +1.37%  "
$ws.Range("D10").Value = "'17.97"
$ws.Range("E10").Value = "  +12.11%  "
$ws.Range("D11").Value = "'1.838.93"
$ws.Range("E11").Value = "  +0.51%  "
$ws.Range("E12").Value = "  +0.70%  "
$ws.Range("D13").Value = "'4.994"
$ws.Range("E13").Value = "  +2.24%  "
$ws.Range("D14").Value = "'84.84"
$ws.Range("E14").Value = "  +2.00%  "
$ws.Range("D15").Value = "'0.6244"
$ws.Range("E15").Value = "  +1.40%  "
$ws.Range("D16").Value = "'30.438.40"
$ws.Range("E16").Value = "  +1.20%  "
$ws.Range("D17").Value = "'252.71"
$ws.Range("E17").Value = "  +11.82%  "
$ws.Range("D18").Value = "'0.9995"
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("D19").Value = "'12.70"
$ws.Range("E19").Value = "  +2.65%  "
$ws.Range("D20").Value = "'0.000007353"
$ws.Range("E20").Value = "  +1.44%  "
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("D22").Value = "'4.934"
$ws.Range("E22").Value = "  +1.92%  "
$ws.Range("D23").Value = "'5.906"
$ws.Range("E23").Value = "  +1.39%  "
$ws.Range("D24").Value = "'164.07"
$ws.Range("E24").Value = "  -0.42%  "
$ws.Range("D25").Value = "'9.001"
$ws.Range("E25").Value = "  -1.70%  "
$ws.Range("D26").Value = "'18.02"
$ws.Range("E26").Value = "  +2.75%  "
$ws.Range("D27").Value = "'1.881"
$ws.Range("E27").Value = "  +1.20%  "
$ws.Range("D28").Value = "'0.1025"
$ws.Range("E28").Value = "  +1.84%  "
$ws.Range("E29").Value = "  -1.66%  "
$ws.Range("D30").Value = "'4.048"
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("D31").Value = "'3.841"
$ws.Range("E31").Value = "  +1.84%  "
$ws.Range("D32").Value = "'0.04834"
$ws.Range("E32").Value = "  +1.46%  "
$ws.Range("D33").Value = "'1.133"
$ws.Range("E33").Value = "  +0.60%  "
$ws.Range("D34").Value = "'0.6992"
$ws.Range("E34").Value = "  -1.19%  "
$ws.Range("D35").Value = "'2.692"
$ws.Range("E35").Value = "  -0.41%  "
$ws.Range("D36").Value = "'0.01888"
$ws.Range("E36").Value = "  +4.49%  "
$ws.Range("D37").Value = "'2.683"
$ws.Range("E37").Value = "  +2.79%  "
$ws.Range("D38").Value = "'0.8748"
$ws.Range("E38").Value = "  -1.72%  "
$ws.Range("D39").Value = "'1.984"
$ws.Range("E39").Value = "  +3.13%  "
$ws.Range("E40").Value = "  +2.74%  "
$ws.Range("D41").Value = "'0.9998"
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").Value = "'0.4072"
$ws.Range("E42").Value = "  +2.10%  "
$ws.Range("D43").Value = "'5.509"
$ws.Range("E43").Value = "  +0.79%  "
$ws.Range("D44").Value = "'7.171"
$ws.Range("E44").Value = "  +3.32%  "
$ws.Range("D45").Value = "'63.20"
$ws.Range("E45").Value = "  +6.10%  "
$ws.Range("E46").Value = "  +1.27%  "
$ws.Range("D47").Value = "'34.11"
$ws.Range("E47").Value = "  +4.38%  "
$ws.Range("D48").Value = "'8.564"
$ws.Range("E48").Value = "  +2.33%  "
$ws.Range("D49").Value = "'0.05502"
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("D50").Value = "'1.352"
$ws.Range("E50").Value = "  -0.95%  "
$ws.Range("D51").Value = "'0.3693"
$ws.Range("E51").Value = "  +2.20%  "
